$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1) - update "想去人数" (column F) counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 118
$wsExhibit.Range("F4").Value = 1628
$wsExhibit.Range("F5").Value = 624
$wsExhibit.Range("F6").Value = 1097
$wsExhibit.Range("F8").Value = 11537
$wsExhibit.Range("F9").Value = 28
$wsExhibit.Range("F11").Value = 452
$wsExhibit.Range("F12").Value = 364
$wsExhibit.Range("F14").Value = 804
$wsExhibit.Range("F15").Value = 12390
$wsExhibit.Range("F16").Value = 13106
$wsExhibit.Range("F21").Value = 230
$wsExhibit.Range("F24").Value = 126

# "全部类型" sheet (sheet4) - same rows/events, same counts (F3 already 118, unchanged)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1628
$wsAll.Range("F5").Value = 624
$wsAll.Range("F6").Value = 1097
$wsAll.Range("F8").Value = 11537
$wsAll.Range("F9").Value = 28
$wsAll.Range("F11").Value = 452
$wsAll.Range("F12").Value = 364
$wsAll.Range("F14").Value = 804
$wsAll.Range("F15").Value = 12390
$wsAll.Range("F16").Value = 13106
$wsAll.Range("F21").Value = 230
$wsAll.Range("F24").Value = 126
